{"js": "// Updates the \"Ergebnisdokumentation zu React\" paragraphs:\n//  1) Eventhandler/States paragraph gets expanded explanation text.\n//  2) The \"React merkt, dass sich ... sichtbar.\" sentence is reworded.\n//  3) The top-down/bottom-up Informationsfluss sentence is reworded\n//     (split into two search/replace calls so the \"bottom-up\" run in\n//     between -- which Word's spell-checker wraps separately -- is left\n//     untouched).\n//  4) The \"Hooks ...\" placeholder paragraph is replaced with a paragraph\n//     containing a Wingdings arrow symbol and a new remark.\n\nconst body = context.document.body;\n\nasync function replaceOnce(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Eventhandler / States paragraph (text up to, but excluding, \"React\").\nawait replaceOnce(\n  \" Eine Ausnahme sind dabei Eventhandler. Da diese nicht w\u00e4hrend des Renderns laufen, m\u00fcssen sie nicht pure sein und k\u00f6nnen deshalb \u00c4nderungen vornehmen. Jedoch ver\u00e4ndern auch sie keine Variablen, sondern sogenannte States. Der State ist eine Art \u00abGed\u00e4chtnis\u00bb einer Komponente. Wenn \",\n  \" Eine Ausnahme sind dabei Eventhandler. Da sie nicht w\u00e4hrend des Renderns laufen, m\u00fcssen sie nicht pure sein und k\u00f6nnen deshalb \u00c4nderungen vornehmen. Damit Daten ge\u00e4ndert werden k\u00f6nnen, m\u00fcssen sie jedoch erst irgendwo gespeichert werden. Daf\u00fcr werden States verwendet. States sind eine Art \u00abGed\u00e4chtnis\u00bb einer Komponente. Wenn \"\n);\n\n// 2) \"React merkt, dass sich ... sofort sichtbar.\" sentence (text after \"React\").\nawait replaceOnce(\n  \" merkt, dass sich etwas ge\u00e4ndert hat (z. B. ein State), wird die ganze Komponente neu gerendert und die \u00c4nderungen werden sofort sichtbar.\",\n  \" merkt, dass sich ein solcher State ver\u00e4ndert hat, wird die ganze Komponente neu gerendert und die \u00c4nderung wird sofort sichtbar.\"\n);\n\n// 3) Top-down Informationsfluss sentence (text up to, but excluding, \"bottom-up\").\nawait replaceOnce(\n  \" und dienen dem top-down Informationsfluss zwischen einzelnen Komponenten. Eine weitere M\u00f6glichkeit f\u00fcr den Informationsfluss zwischen Komponenten findet \",\n  \" und dienen dem top-down Informationsfluss zwischen einzelnen Komponenten. F\u00fcr den \"\n);\n\n// 4) Bottom-up Informationsfluss sentence (text after \"bottom-up\").\nawait replaceOnce(\n  \" mit Hilfe von Callback-Methoden statt\",\n  \" Informationsfluss zwischen Komponenten werden sogenannte Callback-Methoden verwendet\"\n);\n\n// 5) Replace the \"Hooks ...\" placeholder paragraph with the updated\n//    remark paragraph (left-aligned, Wingdings arrow, new comment text).\nconst hooksResults = body.search(\"Hooks\", { matchCase: true });\nhooksResults.load(\"items\");\nawait context.sync();\n\nconst hooksParagraph = hooksResults.items[0].paragraphs.getFirst();\nconst newParagraphXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  '<w:pPr><w:jc w:val=\"left\"/></w:pPr>' +\n  \"<w:r><w:t>Hooks</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"fr-CH\"/></w:rPr><w:sym w:font=\"Wingdings\" w:char=\"F0E0\"/></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> wichtig oder evtl. W</w:t></w:r>' +\n  \"<w:r><w:t>eglassen?</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nhooksParagraph.getRange().insertOoxml(newParagraphXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Updates the \"Ergebnisdokumentation zu React\" paragraphs:\n#  1) Eventhandler/States paragraph gets expanded explanation text.\n#  2) The \"React merkt, dass sich ... sichtbar.\" sentence is reworded.\n#  3) The top-down/bottom-up Informationsfluss sentence is reworded\n#     (split into two Find/Replace calls so the \"bottom-up\" proofErr run\n#     in between is left untouched).\n#  4) The \"Hooks ...\" placeholder paragraph is replaced with a paragraph\n#     containing a Wingdings arrow symbol and a new remark.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# 1) Eventhandler / States paragraph (text up to, but excluding, \"React\").\nReplace-Text `\n    \" Eine Ausnahme sind dabei Eventhandler. Da diese nicht w\u00e4hrend des Renderns laufen, m\u00fcssen sie nicht pure sein und k\u00f6nnen deshalb \u00c4nderungen vornehmen. Jedoch ver\u00e4ndern auch sie keine Variablen, sondern sogenannte States. Der State ist eine Art \u00abGed\u00e4chtnis\u00bb einer Komponente. Wenn \" `\n    \" Eine Ausnahme sind dabei Eventhandler. Da sie nicht w\u00e4hrend des Renderns laufen, m\u00fcssen sie nicht pure sein und k\u00f6nnen deshalb \u00c4nderungen vornehmen. Damit Daten ge\u00e4ndert werden k\u00f6nnen, m\u00fcssen sie jedoch erst irgendwo gespeichert werden. Daf\u00fcr werden States verwendet. States sind eine Art \u00abGed\u00e4chtnis\u00bb einer Komponente. Wenn \"\n\n# 2) \"React merkt, dass sich ... sofort sichtbar.\" sentence (text after \"React\").\nReplace-Text `\n    \" merkt, dass sich etwas ge\u00e4ndert hat (z. B. ein State), wird die ganze Komponente neu gerendert und die \u00c4nderungen werden sofort sichtbar.\" `\n    \" merkt, dass sich ein solcher State ver\u00e4ndert hat, wird die ganze Komponente neu gerendert und die \u00c4nderung wird sofort sichtbar.\"\n\n# 3) Top-down Informationsfluss sentence (text up to, but excluding, \"bottom-up\").\nReplace-Text `\n    \" und dienen dem top-down Informationsfluss zwischen einzelnen Komponenten. Eine weitere M\u00f6glichkeit f\u00fcr den Informationsfluss zwischen Komponenten findet \" `\n    \" und dienen dem top-down Informationsfluss zwischen einzelnen Komponenten. F\u00fcr den \"\n\n# 4) Bottom-up Informationsfluss sentence (text after \"bottom-up\").\nReplace-Text `\n    \" mit Hilfe von Callback-Methoden statt\" `\n    \" Informationsfluss zwischen Komponenten werden sogenannte Callback-Methoden verwendet\"\n\n# 5) Replace the \"Hooks ...\" placeholder paragraph with the updated\n#    remark paragraph (left-aligned, Wingdings arrow, new comment text).\n$hooksFind = $d.Content.Find\n$hooksFind.Text = \"Hooks\"\n$hooksFind.Execute() | Out-Null\n$hooksPara = $hooksFind.Parent.Paragraphs(1).Range\n\n$newParaXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:jc w:val=\"left\"/></w:pPr><w:r><w:t>Hooks</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:lang w:val=\"fr-CH\"/></w:rPr><w:sym w:font=\"Wingdings\" w:char=\"F0E0\"/></w:r><w:r><w:t xml:space=\"preserve\"> wichtig oder evtl. W</w:t></w:r><w:r><w:t>eglassen?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$hooksPara.InsertXML($newParaXml)\n"}
